$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Address, $Value)
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

$updates = [ordered]@{
    'D2' = '29.971.39'
    'E2' = '  -1.00%  '
    'D3' = '1.878.25'
    'E3' = '  -2.01%  '
    'D4' = '0.9986'
    'E4' = '  -0.03%  '
    'D5' = '242.32'
    'E5' = '  -4.55%  '
    'D6' = '0.9982'
    'E6' = '  -0.07%  '
    'D7' = '0.4985'
    'E7' = '  -3.12%  '
    'E8' = '  -2.02%  '
    'D9' = '0.06630'
    'E9' = '  -2.96%  '
    'D10' = '1.877.35'
    'E10' = '  -1.98%  '
    'D11' = '16.73'
    'E11' = '  -4.19%  '
    'D12' = '0.07244'
    'E12' = '  -1.58%  '
    'D13' = '0.6681'
    'E13' = '  -3.99%  '
    'D14' = '86.30'
    'E14' = '  -1.79%  '
    'D15' = '4.884'
    'E15' = '  -0.48%  '
    'D16' = '29.945.42'
    'E16' = '  -1.06%  '
    'D17' = '0.000007926'
    'E17' = '  -0.63%  '
    'E18' = '  -0.13%  '
    'E19' = '  -2.24%  '
    'D20' = '2.120.17'
    'E20' = '  -1.95%  '
    'D21' = '0.9983'
    'E21' = '  -0.01%  '
    'D22' = '4.762'
    'E22' = '  -2.08%  '
    'D23' = '5.654'
    'E23' = '  -1.56%  '
    'E24' = '  -1.36%  '
    'D25' = '149.16'
    'E25' = '  +1.89%  '
    'D26' = '141.92'
    'E26' = '  +2.18%  '
    'D27' = '17.21'
    'E27' = '  -0.61%  '
    'D28' = '1.913'
    'E28' = '  -5.61%  '
    'D29' = '1.389'
    'E29' = '  +0.28%  '
    'D30' = '4.172'
    'E30' = '  -2.55%  '
    'D31' = '0.08778'
    'E31' = '  -0.77%  '
    'D32' = '3.943'
    'E32' = '  -2.24%  '
    'D33' = '0.05074'
    'E33' = '  -1.23%  '
    'D34' = '0.7098'
    'E34' = '  -1.25%  '
    'E35' = '  -4.73%  '
    'D36' = '2.664'
    'B37' = 'VeChain'
    'C37' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D37' = '0.01756'
    'E37' = '  +3.20%  '
    'B38' = 'MXToken'
    'C38' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D38' = '2.689'
    'E38' = '  -5.49%  '
    'D39' = '2.180'
    'E39' = '  -5.97%  '
    'E40' = '  -4.85%  '
    'B41' = 'FraxShare'
    'C41' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D41' = '5.791'
    'E41' = '  -5.17%  '
    'B42' = 'TheSandbox'
    'C42' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D42' = '0.4252'
    'E42' = '  -1.79%  '
    'E43' = '  -0.13%  '
    'D44' = '102.14'
    'E44' = '  -4.03%  '
    'D45' = '7.467'
    'E45' = '  -3.56%  '
    'D46' = '0.1258'
    'E46' = '  -2.18%  '
    'D47' = '0.05652'
    'E47' = '  -1.83%  '
    'E48' = '  -3.32%  '
    'D49' = '0.3761'
    'E49' = '  -2.19%  '
    'D50' = '8.199'
    'E50' = '  -4.42%  '
    'D51' = '55.84'
    'E51' = '  -2.00%  '
}

foreach ($addr in $updates.Keys) {
    Set-TextValue $ws $addr $updates[$addr]
}
